$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37691
$ws.Range("D2").Value = 54510576
$ws.Range("C3").Value = 90890
$ws.Range("D3").Value = 133238139
$ws.Range("C4").Value = 31151
$ws.Range("D4").Value = 46133935
$ws.Range("C5").Value = 8684
$ws.Range("D5").Value = 12906563
$ws.Range("C6").Value = 1990
$ws.Range("D6").Value = 2957506
$ws.Range("C12").Value = 41276
$ws.Range("D12").Value = 56003689
$ws.Range("C13").Value = 9647
$ws.Range("D13").Value = 13952458
$ws.Range("C14").Value = 25929
$ws.Range("D14").Value = 38029540
$ws.Range("C15").Value = 8307
$ws.Range("D15").Value = 12328324
$ws.Range("C16").Value = 2150
$ws.Range("D16").Value = 3197165
$ws.Range("C20").Value = 10209
$ws.Range("D20").Value = 13520278
$ws.Range("C21").Value = 13370
$ws.Range("D21").Value = 19305292
$ws.Range("C22").Value = 31632
$ws.Range("D22").Value = 46421595
$ws.Range("C23").Value = 10212
$ws.Range("D23").Value = 15181055
$ws.Range("C27").Value = 11673
$ws.Range("D27").Value = 15593409
$ws.Range("C28").Value = 7634
$ws.Range("D28").Value = 11058377
$ws.Range("C29").Value = 22463
$ws.Range("D29").Value = 32972549
$ws.Range("C30").Value = 7807
$ws.Range("D30").Value = 11618133
$ws.Range("C34").Value = 8300
$ws.Range("D34").Value = 10963515
$ws.Range("C35").Value = 3239
$ws.Range("D35").Value = 4675194
$ws.Range("C36").Value = 7819
$ws.Range("D36").Value = 11418956
$ws.Range("C37").Value = 3175
$ws.Range("D37").Value = 4705461
$ws.Range("C41").Value = 2471
$ws.Range("D41").Value = 3339853
$ws.Range("C42").Value = 17217
$ws.Range("D42").Value = 24897303
$ws.Range("C43").Value = 51046
$ws.Range("D43").Value = 74833767
$ws.Range("C44").Value = 19000
$ws.Range("D44").Value = 28222443
$ws.Range("C45").Value = 5603
$ws.Range("D45").Value = 8343677
$ws.Range("C46").Value = 1201
$ws.Range("D46").Value = 1792045
$ws.Range("C50").Value = 16674
$ws.Range("D50").Value = 22200892
$ws.Range("C51").Value = 2013
$ws.Range("D51").Value = 2919918
$ws.Range("C52").Value = 6877
$ws.Range("D52").Value = 10109074
$ws.Range("C53").Value = 2344
$ws.Range("D53").Value = 3500918
$ws.Range("C54").Value = 754
$ws.Range("D54").Value = 1126305
$ws.Range("C55").Value = 185
$ws.Range("D55").Value = 274333
$ws.Range("C57").Value = 6936
$ws.Range("D57").Value = 9539107
$ws.Range("C58").Value = 935
$ws.Range("D58").Value = 1372079
$ws.Range("C59").Value = 2362
$ws.Range("D59").Value = 3501837
$ws.Range("C60").Value = 938
$ws.Range("D60").Value = 1396501
$ws.Range("C64").Value = 1384
$ws.Range("D64").Value = 1947206
$ws.Range("C65").Value = 15331
$ws.Range("D65").Value = 22146743
$ws.Range("C66").Value = 44635
$ws.Range("D66").Value = 65317207
$ws.Range("C67").Value = 15686
$ws.Range("D67").Value = 23311186
$ws.Range("C68").Value = 4564
$ws.Range("D68").Value = 6798292
$ws.Range("C69").Value = 919
$ws.Range("D69").Value = 1366668
$ws.Range("C73").Value = 15065
$ws.Range("D73").Value = 19865494
$ws.Range("C74").Value = 51274
$ws.Range("D74").Value = 74615550
$ws.Range("C75").Value = 145770
$ws.Range("D75").Value = 214755690
$ws.Range("C76").Value = 63529
$ws.Range("D76").Value = 94667865
$ws.Range("C77").Value = 20306
$ws.Range("D77").Value = 30339331
$ws.Range("C78").Value = 4807
$ws.Range("D78").Value = 7179543
$ws.Range("C85").Value = 50712
$ws.Range("D85").Value = 68986653
$ws.Range("C86").Value = 4590
$ws.Range("D86").Value = 6650436
$ws.Range("C87").Value = 11545
$ws.Range("D87").Value = 16961042
$ws.Range("C88").Value = 3881
$ws.Range("D88").Value = 5784083
$ws.Range("C93").Value = 5400
$ws.Range("D93").Value = 7259391
$ws.Range("C94").Value = 1591
$ws.Range("D94").Value = 2291432
$ws.Range("C95").Value = 5146
$ws.Range("D95").Value = 7577743
$ws.Range("C96").Value = 1938
$ws.Range("D96").Value = 2886937
$ws.Range("C101").Value = 3543
$ws.Range("D101").Value = 4687961
$ws.Range("C102").Value = 600
$ws.Range("D102").Value = 893664
$ws.Range("C103").Value = 350
$ws.Range("D103").Value = 522530
$ws.Range("C105").Value = 44
$ws.Range("D105").Value = 66000
$ws.Range("C107").Value = 10738
$ws.Range("D107").Value = 15576962
$ws.Range("C108").Value = 29164
$ws.Range("D108").Value = 42850642
$ws.Range("C109").Value = 9766
$ws.Range("D109").Value = 14522650
$ws.Range("C110").Value = 2679
$ws.Range("D110").Value = 3994707
$ws.Range("C114").Value = 9779
$ws.Range("D114").Value = 12920313
$ws.Range("C115").Value = 30379
$ws.Range("D115").Value = 43810202
$ws.Range("C116").Value = 66070
$ws.Range("D116").Value = 96695393
$ws.Range("C117").Value = 21345
$ws.Range("D117").Value = 31721940
$ws.Range("C118").Value = 6057
$ws.Range("D118").Value = 9023521
$ws.Range("C124").Value = 25802
$ws.Range("D124").Value = 34465594
$ws.Range("C125").Value = 35924
$ws.Range("D125").Value = 51848102
$ws.Range("C126").Value = 76678
$ws.Range("D126").Value = 112127243
$ws.Range("C127").Value = 23812
$ws.Range("D127").Value = 35339409
$ws.Range("C128").Value = 6384
$ws.Range("D128").Value = 9486738
$ws.Range("C130").Value = 58
$ws.Range("D130").Value = 85228
$ws.Range("C133").Value = 31748
$ws.Range("D133").Value = 42161709
$ws.Range("C134").Value = 13195
$ws.Range("D134").Value = 19099407
$ws.Range("C135").Value = 32280
$ws.Range("D135").Value = 47413893
$ws.Range("C136").Value = 11459
$ws.Range("D136").Value = 17026542
$ws.Range("C137").Value = 2953
$ws.Range("D137").Value = 4402214
$ws.Range("C138").Value = 498
$ws.Range("D138").Value = 740990
$ws.Range("C141").Value = 10799
$ws.Range("D141").Value = 14402256
$ws.Range("C142").Value = 34964
$ws.Range("D142").Value = 50491751
$ws.Range("C143").Value = 81100
$ws.Range("D143").Value = 118824486
$ws.Range("C144").Value = 24305
$ws.Range("D144").Value = 36112055
$ws.Range("C145").Value = 6379
$ws.Range("D145").Value = 9518067
$ws.Range("C146").Value = 1430
$ws.Range("D146").Value = 2127230
$ws.Range("C149").Value = 29130
$ws.Range("D149").Value = 39303536
